$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update the borrower name and the headline figures.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = "Kamal Al Suwaidi"
$wsSummary.Range("B4").Value = 2234.8
$wsSummary.Range("B6").Value = 935224
$wsSummary.Range("B7").Value = 258170
$wsSummary.Range("B8").Value = 677054
$wsSummary.Range("B9").Value = 3.62

# ---------------------------------------------------------------------------
# Sheet "Assets": insert two new "Vehicles / Luxury Car" rows above the
# existing "Liquid Assets" row, then update the values that shifted down.
# ---------------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")

# Push the existing data rows (Liquid Assets, TOTAL ASSETS) down by two rows.
$wsAssets.Rows.Item(2).Insert()
$wsAssets.Rows.Item(2).Insert()

# New row 2: Vehicles / Luxury Car
$wsAssets.Range("A2").Value = "Vehicles"
$wsAssets.Range("B2").Value = "Luxury Car"
$wsAssets.Range("C2").Value = 365578

# New row 3: Vehicles / Luxury Car
$wsAssets.Range("A3").Value = "Vehicles"
$wsAssets.Range("B3").Value = "Luxury Car"
$wsAssets.Range("C3").Value = 566613

# The new rows inherit a generic style from the Insert(); copy the real
# formatting (borders/fill/font/number format) from the row below - which
# kept the original "Liquid Assets" row formatting - onto the new rows.
$wsAssets.Range("A4:C4").Copy()
$wsAssets.Range("A2:C2").PasteSpecial(-4122)
$wsAssets.Range("A3:C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 4 (previously row 2, "Liquid Assets"): update the value.
$wsAssets.Range("C4").Value = 3033

# Row 5 (previously row 3, "TOTAL ASSETS"): update the total.
$wsAssets.Range("C5").Value = 935224

# ---------------------------------------------------------------------------
# Sheet "Liabilities": first loan becomes an auto loan, amounts updated.
# ---------------------------------------------------------------------------
$wsLiabilities = $wb.Worksheets.Item("Liabilities")
$wsLiabilities.Range("A2").Value = "Auto Loans"
$wsLiabilities.Range("B2").Value = "Vehicle Loan 1"
$wsLiabilities.Range("C2").Value = 219347
$wsLiabilities.Range("D2").Value = 3656

$wsLiabilities.Range("C3").Value = 38823
$wsLiabilities.Range("D3").Value = 1941

$wsLiabilities.Range("C4").Value = 258170
